$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Complete?" column (A) as done ("x") for the first batch of
# short-effort issues that have now been finished (rows 2-8).
$ws.Range("A2:A8").Value = "x"

# Move the active selection to A12 to reflect where work continues.
$ws.Range("A12").Select()
